$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (columns mirror the sheet header:
# A=ملاحظات, B=المرافق, C=الكمية, D=المخيم, E=نوع المسافة, F=المركبة, G=المؤسسة, H=الوقت)
$newRows = @(
    @("", "احمد", "23", "الصمود", "الرحلة 3", "C3", "NRC", "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٥٥:١٤ م"),
    @("", "احمد", "23", "الصمود", "الرحلة 3", "C3", "NRC", "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٥٥:١٤ م"),
    @("23", "احمد", "2323", "الصمود", "الرحلة 3", "C3", "NRC", "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٥٥:١٤ م"),
    @("", "احمد", "34", "الصمود", "الرحلة 3", "C3", "NRC", "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٥٥:١٤ م")
)

$startRow = 17
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $col = $j + 1
        $cell = $ws.Cells.Item($r, $col)
        $v = $rowValues[$j]
        # Every value in this sheet (including blanks and digit-only
        # quantities such as "23") is stored as literal text rather than a
        # number, so force a text number format before writing anything
        # that Excel would otherwise auto-convert to a numeric value.
        if ($v -eq "" -or $v -match '^[0-9]+$') {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $v
    }
}
